# ----------------------------------------------------------------------
# Add a new worksheet "ODI Batting Extra" as the 3rd sheet (after the
# existing "Player Info" / "ODI Batting" sheets) and populate it with
# per-match batting detail data.
# ----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# Grab a reference header cell (style carries bold font + border used
# for every header row in this workbook) before we add/alter anything.
$headerStyleSource = $wb.Worksheets.Item(1).Range("A1")

# Copy the last existing sheet so the new sheet inherits the same
# sheetPr / pageMargins / sheetFormatPr structure used elsewhere in the
# workbook, then drop straight into the copy and wipe its contents.
$lastIndex = $wb.Worksheets.Count
$wb.Worksheets.Item($lastIndex).Copy([Type]::Missing, $wb.Worksheets.Item($lastIndex))
$ws = $wb.Worksheets.Item($lastIndex + 1)
$ws.Name = "ODI Batting Extra"
$ws.Cells.Clear()

# ---- Header row -------------------------------------------------------
$headers = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($c = 1; $c -le $headers.Length; $c++) {
    $ws.Cells.Item(1, $c).Value = $headers[$c - 1]
}
$headerRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item(1, $headers.Length))
$headerStyleSource.Copy()
$headerRange.PasteSpecial(-4122)   # xlPasteFormats

# ---- Data rows ----------------------------------------------------------
# Columns: MATCH_CODE, BATTING_POSITION, NUM_4, NUM_6, PERCENT_RUNS_OF_TOTAL, MAN_OF_MATCH
# MATCH_CODE / NUM_4 / NUM_6 / PERCENT_RUNS_OF_TOTAL are stored as literal
# text (numeric-looking strings need a leading apostrophe so they are not
# reinterpreted as numbers); BATTING_POSITION is a genuine number;
# MAN_OF_MATCH is plain text.
$data = @(
    @("3829", 1,    "7", "1", "31.08%", "NO"),
    @("3830", 1,    "5", "0", "14.43%", "NO"),
    @("3831", 1,    "4", "0", "7.12%",  "NO"),
    @("3832", $null, $null, $null, $null, "NO"),
    @("3833", 1,    "0", "0", "0.67%",  "NO"),
    @("3834", 1,    "0", "0", $null,    "NO")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowNum = $r + 2
    $row = $data[$r]

    # A: MATCH_CODE - literal text
    $ws.Cells.Item($rowNum, 1).Value = "'" + $row[0]

    # B: BATTING_POSITION - number, blank on row 5 (match 3832)
    if ($row[1] -ne $null) {
        $ws.Cells.Item($rowNum, 2).Value = $row[1]
    } else {
        $ws.Cells.Item($rowNum, 2).Value = "'"
    }

    # C: NUM_4 - literal text (blank on row 5)
    if ($row[2] -ne $null) {
        $ws.Cells.Item($rowNum, 3).Value = "'" + $row[2]
    } else {
        $ws.Cells.Item($rowNum, 3).Value = "'"
    }

    # D: NUM_6 - literal text (blank on row 5)
    if ($row[3] -ne $null) {
        $ws.Cells.Item($rowNum, 4).Value = "'" + $row[3]
    } else {
        $ws.Cells.Item($rowNum, 4).Value = "'"
    }

    # E: PERCENT_RUNS_OF_TOTAL - literal text (blank on rows 5 & 7)
    if ($row[4] -ne $null) {
        $ws.Cells.Item($rowNum, 5).Value = "'" + $row[4]
    } else {
        $ws.Cells.Item($rowNum, 5).Value = "'"
    }

    # F: MAN_OF_MATCH - plain text, always populated
    $ws.Cells.Item($rowNum, 6).Value = $row[5]
}
